# Apply the "finished determining MMOI experimentally" edit to the Notes sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notes")
$ws.Activate()

# --- A11 gets the same center/center alignment already used by the other
#     Topic cells in column A (A2..A10); copy that formatting over. ---
$ws.Range("A2").Copy()
$ws.Range("A11").PasteSpecial(-4122)

# --- New row 12: merge A11:A12 (continuing the Topic cell downward),
#     write the date and the new note text. ---
$ws.Range("A11:A12").Merge()

$ws.Range("B9").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("B12").Value = 45829

$ws.Range("C12").Value = "I tried to replicate the experiment using the LED light switch. However, the error between the experiemnt MMOI and analytical MMOI was larger from my experiment than the error presented by the experiement done in the paper. I think the errors are larger because this remote controller's mass is not evenely distributed like a simple block, and for the analytical method I used a simple block formula to determine the analytical MMOI"

$ws.Rows.Item(12).RowHeight = 90

# --- Match the final selection recorded in the sheet view. ---
$ws.Range("B12").Select()
